$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($range, $value) {
    $range.NumberFormat = "@"
    $range.Value = $value
    $range.Style = "Normal"
}

Set-TextValue $ws.Range("D2") '61.568.73'
Set-TextValue $ws.Range("E2") '  +0.96%  '

Set-TextValue $ws.Range("D3") '3.389.00'
Set-TextValue $ws.Range("E3") '  -0.25%  '

Set-TextValue $ws.Range("D4") '1.00'
Set-TextValue $ws.Range("E4") '  +0.05%  '

Set-TextValue $ws.Range("D5") '575.65'
Set-TextValue $ws.Range("E5") '  +0.47%  '

Set-TextValue $ws.Range("D6") '140.75'
Set-TextValue $ws.Range("E6") '  -1.51%  '

Set-TextValue $ws.Range("E7") '  +0.03%  '

Set-TextValue $ws.Range("E8") '  -0.48%  '

Set-TextValue $ws.Range("D9") '7.72'
Set-TextValue $ws.Range("E9") '  +2.15%  '

Set-TextValue $ws.Range("E10") '  -1.08%  '

Set-TextValue $ws.Range("E11") '  -2.21%  '

Set-TextValue $ws.Range("D12") '3.972.02'
Set-TextValue $ws.Range("E12") '  -0.02%  '

Set-TextValue $ws.Range("D13") '28.48'
Set-TextValue $ws.Range("E13") '  +0.96%  '

Set-TextValue $ws.Range("E14") '  +0.19%  '

Set-TextValue $ws.Range("D15") '3.375.45'
Set-TextValue $ws.Range("E15") '  -0.72%  '

Set-TextValue $ws.Range("E16") '  -0.65%  '

Set-TextValue $ws.Range("D17") '61.575.93'
Set-TextValue $ws.Range("E17") '  +0.89%  '

Set-TextValue $ws.Range("E18") '  -0.43%  '

Set-TextValue $ws.Range("D19") '13.64'
Set-TextValue $ws.Range("E19") '  -1.59%  '

Set-TextValue $ws.Range("D20") '9.00'
Set-TextValue $ws.Range("E20") '  +0.17%  '

Set-TextValue $ws.Range("D21") '391.95'
Set-TextValue $ws.Range("E21") '  +2.27%  '

Set-TextValue $ws.Range("D22") '75.28'
Set-TextValue $ws.Range("E22") '  +1.40%  '

Set-TextValue $ws.Range("D23") '0.554'
Set-TextValue $ws.Range("E23") '  -1.02%  '

Set-TextValue $ws.Range("E24") '  -0.04%  '

Set-TextValue $ws.Range("E25") '  -5.24%  '

Set-TextValue $ws.Range("E26") '  +7.69%  '

Set-TextValue $ws.Range("D27") '0.998'
Set-TextValue $ws.Range("E27") '  -0.11%  '

Set-TextValue $ws.Range("D28") '7.26'
Set-TextValue $ws.Range("E28") '  -1.97%  '

Set-TextValue $ws.Range("D29") '8.06'
Set-TextValue $ws.Range("E29") '  +0.74%  '

Set-TextValue $ws.Range("D30") '2.15'
Set-TextValue $ws.Range("E30") '  -0.59%  '

Set-TextValue $ws.Range("D31") '1.40'
Set-TextValue $ws.Range("E31") '  -1.44%  '

Set-TextValue $ws.Range("E32") '  -0.05%  '

Set-TextValue $ws.Range("E33") '  -1.02%  '

Set-TextValue $ws.Range("E34") '  -1.26%  '

Set-TextValue $ws.Range("D35") '168.26'
Set-TextValue $ws.Range("E35") '  +0.18%  '

Set-TextValue $ws.Range("D36") '5.06'
Set-TextValue $ws.Range("E36") '  +0.99%  '

Set-TextValue $ws.Range("D37") '3.423.90'

Set-TextValue $ws.Range("E38") '  -1.48%  '

Set-TextValue $ws.Range("D39") '0.0768'
Set-TextValue $ws.Range("E39") '  -0.72%  '

Set-TextValue $ws.Range("D40") '26.12'
Set-TextValue $ws.Range("E40") '  -6.01%  '

Set-TextValue $ws.Range("E41") '  -0.33%  '

Set-TextValue $ws.Range("E42") '  -0.47%  '

Set-TextValue $ws.Range("E43") '  -1.28%  '

Set-TextValue $ws.Range("E44") '  +2.00%  '

Set-TextValue $ws.Range("D45") '2.457.06'
Set-TextValue $ws.Range("E45") '  -0.95%  '

Set-TextValue $ws.Range("E46") '  -0.04%  '

Set-TextValue $ws.Range("D47") '6.67'
Set-TextValue $ws.Range("E47") '  -2.32%  '

Set-TextValue $ws.Range("E48") '  +0.07%  '

Set-TextValue $ws.Range("E49") '  -1.58%  '

Set-TextValue $ws.Range("E50") '  -4.67%  '

Set-TextValue $ws.Range("E51") '  -1.82%  '
